# Apply the "Official Proceeding" edit: update company/platoon/squad
# numbers for the first 7 roster rows, fill in the commander + deputy
# rows (7-8) with their post, rank, name and DOB, drop the now-unused
# trailing rows 9-11, unhide column B, and reset the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Отделение (squad) changes from 1 to 3 ---
$ws.Range("E2").Value = 3

# --- Row 3: Взвод (platoon) changes from 3 to 1 ---
$ws.Range("D3").Value = 1

# --- Row 4: Взвод (platoon) changes from 4 to 1 ---
$ws.Range("D4").Value = 1

# --- Row 6: Рота (company) changes from 2 to 3 ---
$ws.Range("C6").Value = 3

# --- Row 7: becomes the Company Commander row ---
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2

# --- Row 8: becomes the Deputy Commander (political) row ---
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 1

# Fill in the new personnel details column-by-column (names, then
# positions, then birth dates) to mirror the authoring order.
$ws.Range("N7").Value = "Даутов Искандер Садыкович"
$ws.Range("N8").Value = "Говоров Леонид Александрович"

$ws.Range("F7").Value = "Командир роты"
$ws.Range("F8").Value = "заместитель командира роты по военно-политической работе"

$ws.Range("M7").Value = "рядовой"
$ws.Range("M8").Value = "рядовой"

# "06.04.1967" parses as a valid date (6-Apr-1967), so force text via a
# temporary Text number format, then restore the default "Normal" style
# so the saved cell carries no explicit style (matching plain data rows).
$o7 = $ws.Range("O7")
$o7.NumberFormat = "@"
$o7.Value = "06.04.1967"
$o7.Style = "Normal"

$ws.Range("O8").Value = "31.03.1971"

# --- Rows 9-11 no longer used: delete them entirely ---
$ws.Range("A9:AZ11").EntireRow.Delete()

# --- Column B: unhide and give it a real width ---
# Target stored width is 17.73046875 chars; ColumnWidth assignment snaps
# to this host's pixel grid, so feed it the input that lands closest to
# that stored value (grid step here resolves to 17.6666... chars).
$colB = $ws.Range("B1").EntireColumn
$colB.Hidden = $false
$colB.ColumnWidth = 16.8

# --- Reset the active selection to C1 ---
$ws.Range("C1").Select()
